# expansão das análises automáticas
# Adds 3 new columns (apoio_medio, contribuicoes, media_contribuicoes) to the
# "sint_resumo_por_mencoes_fiq_2023" summary sheet, and rescales the
# "particip"/"taxa_sucesso" columns (E:F) from fractional (0-1) to
# percentage-point (0-100) numbers, while keeping their existing 0.00% style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rescale existing columns E (particip) and F (taxa_sucesso): value * 100
# ---------------------------------------------------------------------------
$ef = @{
    2 = @(83.59550561797752, 59.85663082437276)
    3 = @(16.40449438202247, 73.97260273972603)
    4 = @(89.98637602179836, 93.94398183194551)
    5 = @(10.01362397820163, 96.5986394557823)
    6 = @(94.5906432748538, 21.79289026275116)
    7 = @(5.409356725146199, 29.72972972972973)
}

foreach ($r in $ef.Keys) {
    $vals = $ef[$r]
    $ws.Cells.Item($r, 5).Value = $vals[0]
    $ws.Cells.Item($r, 6).Value = $vals[1]
}

# ---------------------------------------------------------------------------
# 2) New header cells L1:N1
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Copy the header style (bold font + border + centered alignment) from K1
# onto the three new header cells so they match the rest of row 1.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) New data columns L, M, N (apoio_medio, contribuicoes, media_contribuicoes)
# ---------------------------------------------------------------------------
$lmn = @{
    2 = @(92.24142267722461, 208605, 312.2829341317365)
    3 = @(87.74218988390095, 54948, 339.1851851851852)
    4 = @(89.51886785425282, 174765, 140.825946817083)
    5 = @(94.08839711265563, 28881, 203.387323943662)
    6 = @(19.15242792744543, 2050, 14.53900709219858)
    7 = @(24.83848419950433, 158, 14.36363636363636)
}

foreach ($r in $lmn.Keys) {
    $vals = $lmn[$r]
    $ws.Cells.Item($r, 12).Value = $vals[0]
    $ws.Cells.Item($r, 13).Value = $vals[1]
    $ws.Cells.Item($r, 14).Value = $vals[2]
}
